# Include K-Factor in BOM
#
# 1. Add a new "K-Factor" header in column G of the Parts sheet.
# 2. Turn on AutoFilter for the header row (A1:F1), which also registers
#    the hidden, sheet-scoped _xlnm._FilterDatabase defined name that
#    Excel writes whenever a filter is applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM column header -> appended to shared strings, dimension/row spans
# grow from A1:F1 to A1:G1 automatically.
$ws.Range("G1").Value = "K-Factor"

# Apply the AutoFilter over the original header range A1:F1 (new column G
# is intentionally left outside the filter range, matching the source
# workbook).
$ws.Range("A1:F1").AutoFilter()

# Excel records the filter range as a hidden, sheet-scoped workbook-level
# name. Mirror that explicitly so it round-trips into <definedNames>.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Parts!`$A`$1:`$F`$1")
$filterName.Visible = $false
